$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 168163.83
$ws.Cells.Item(113, 9).Value = 250750.75
$ws.Cells.Item(113, 10).Value = 2990
$ws.Cells.Item(113, 11).Value = 250750.75
$ws.Cells.Item(113, 12).Value = 2990
$ws.Cells.Item(113, 13).Value = -247496.75
$ws.Cells.Item(113, 14).Value = -9498

$ws.Cells.Item(132, 8).Value = 276789.1
$ws.Cells.Item(132, 9).Value = 311833.56
$ws.Cells.Item(132, 10).Value = 49000
$ws.Cells.Item(132, 11).Value = 935500.6799999999
$ws.Cells.Item(132, 12).Value = 147000
$ws.Cells.Item(132, 13).Value = -932970.6799999999
$ws.Cells.Item(132, 14).Value = -152060

$ws.Cells.Item(135, 8).Value = 8821.714
$ws.Cells.Item(135, 9).Value = 9346.462
$ws.Cells.Item(135, 10).Value = 2000
$ws.Cells.Item(135, 11).Value = 84118.158
$ws.Cells.Item(135, 12).Value = 18000
$ws.Cells.Item(135, 13).Value = -81583.158
$ws.Cells.Item(135, 14).Value = -23070

$ws.Cells.Item(137, 8).Value = 2171.7144
$ws.Cells.Item(137, 9).Value = 1351
$ws.Cells.Item(137, 10).Value = 2500
$ws.Cells.Item(137, 11).Value = 4053
$ws.Cells.Item(137, 12).Value = 7500
$ws.Cells.Item(137, 13).Value = -1503
$ws.Cells.Item(137, 14).Value = -12600

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 4915.931
$ws.Cells.Item(2, 9).Value = 6420.381
$ws.Cells.Item(2, 10).Value = 966.75
$ws.Cells.Item(2, 11).Value = 6420.381
$ws.Cells.Item(2, 12).Value = 966.75
$ws.Cells.Item(2, 13).Value = -6307.381
$ws.Cells.Item(2, 14).Value = -1192.75

$ws.Cells.Item(116, 8).Value = 4915.931
$ws.Cells.Item(116, 9).Value = 6420.381
$ws.Cells.Item(116, 10).Value = 966.75
$ws.Cells.Item(116, 11).Value = 6420.381
$ws.Cells.Item(116, 12).Value = 966.75
$ws.Cells.Item(116, 13).Value = -4126.381
$ws.Cells.Item(116, 14).Value = -5554.75

$ws.Cells.Item(132, 8).Value = 3149.8333
$ws.Cells.Item(132, 9).Value = 2453.8333
$ws.Cells.Item(132, 10).Value = 4193.8335
$ws.Cells.Item(132, 11).Value = 7361.499899999999
$ws.Cells.Item(132, 12).Value = 12581.5005
$ws.Cells.Item(132, 13).Value = -4831.499899999999
$ws.Cells.Item(132, 14).Value = -17641.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 4915.931
$ws.Cells.Item(3, 9).Value = 6420.381
$ws.Cells.Item(3, 10).Value = 966.75
$ws.Cells.Item(3, 11).Value = 6420.381
$ws.Cells.Item(3, 12).Value = 966.75
$ws.Cells.Item(3, 13).Value = -6306.381
$ws.Cells.Item(3, 14).Value = -1194.75

$ws.Cells.Item(64, 8).Value = 459.81818
$ws.Cells.Item(64, 9).Value = 465.75
$ws.Cells.Item(64, 11).Value = 465.75
$ws.Cells.Item(64, 13).Value = -240.75

$ws.Cells.Item(67, 8).Value = 459.81818
$ws.Cells.Item(67, 9).Value = 465.75
$ws.Cells.Item(67, 11).Value = 465.75
$ws.Cells.Item(67, 13).Value = 314.25

$ws.Cells.Item(134, 8).Value = 4601.8
$ws.Cells.Item(134, 9).Value = 2866.3333
$ws.Cells.Item(134, 10).Value = 5345.5713
$ws.Cells.Item(134, 11).Value = 8598.999899999999
$ws.Cells.Item(134, 12).Value = 16036.7139
$ws.Cells.Item(134, 13).Value = -6063.999899999999
$ws.Cells.Item(134, 14).Value = -21106.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 741.1
$ws.Cells.Item(16, 10).Value = 1200
$ws.Cells.Item(16, 12).Value = 1200
$ws.Cells.Item(16, 14).Value = -1774

$ws.Cells.Item(31, 8).Value = 2142.2812
$ws.Cells.Item(31, 9).Value = 1565.7778
$ws.Cells.Item(31, 10).Value = 2367.8696
$ws.Cells.Item(31, 11).Value = 1565.7778
$ws.Cells.Item(31, 12).Value = 2367.8696
$ws.Cells.Item(31, 13).Value = -1270.7778
$ws.Cells.Item(31, 14).Value = -2957.8696

$ws.Cells.Item(32, 8).Value = 50000
$ws.Cells.Item(32, 9).Value = 50000
$ws.Cells.Item(32, 11).Value = 50000
$ws.Cells.Item(32, 13).Value = -49684

$ws.Cells.Item(34, 8).Value = 2142.2812
$ws.Cells.Item(34, 9).Value = 1565.7778
$ws.Cells.Item(34, 10).Value = 2367.8696
$ws.Cells.Item(34, 11).Value = 1565.7778
$ws.Cells.Item(34, 12).Value = 2367.8696
$ws.Cells.Item(34, 13).Value = -1363.7778
$ws.Cells.Item(34, 14).Value = -2771.8696

$ws.Cells.Item(96, 8).Value = 16119.167
$ws.Cells.Item(96, 10).Value = 16119.167
$ws.Cells.Item(96, 12).Value = 16119.167
$ws.Cells.Item(96, 14).Value = -21611.167

$ws.Cells.Item(99, 8).Value = 5684546.5
$ws.Cells.Item(99, 10).Value = 1600
$ws.Cells.Item(99, 12).Value = 1600
$ws.Cells.Item(99, 14).Value = -4596

$ws.Cells.Item(113, 8).Value = 741.1
$ws.Cells.Item(113, 10).Value = 1200
$ws.Cells.Item(113, 12).Value = 1200
$ws.Cells.Item(113, 14).Value = -5540

$ws.Cells.Item(122, 8).Value = 1179.7
$ws.Cells.Item(122, 9).Value = 1149.8334
$ws.Cells.Item(122, 11).Value = 3449.5002
$ws.Cells.Item(122, 13).Value = -999.5002

$ws.Cells.Item(126, 8).Value = 5684546.5
$ws.Cells.Item(126, 10).Value = 1600
$ws.Cells.Item(126, 12).Value = 4800
$ws.Cells.Item(126, 14).Value = -9740

$ws.Cells.Item(134, 8).Value = 2419.7827
$ws.Cells.Item(134, 9).Value = 1022.94446
$ws.Cells.Item(134, 10).Value = 7448.4
$ws.Cells.Item(134, 11).Value = 3068.83338
$ws.Cells.Item(134, 12).Value = 22345.2
$ws.Cells.Item(134, 13).Value = -533.83338
$ws.Cells.Item(134, 14).Value = -27415.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1047.7059
$ws.Cells.Item(5, 9).Value = 482.64285
$ws.Cells.Item(5, 10).Value = 1443.25
$ws.Cells.Item(5, 11).Value = 1447.92855
$ws.Cells.Item(5, 12).Value = 4329.75
$ws.Cells.Item(5, 13).Value = -1335.92855
$ws.Cells.Item(5, 14).Value = -4553.75

$ws.Cells.Item(92, 8).Value = 525
$ws.Cells.Item(92, 9).Value = 200
$ws.Cells.Item(92, 10).Value = 850
$ws.Cells.Item(92, 11).Value = 600
$ws.Cells.Item(92, 12).Value = 2550
$ws.Cells.Item(92, 13).Value = 648
$ws.Cells.Item(92, 14).Value = -5046

$ws.Cells.Item(122, 8).Value = 923.05884
$ws.Cells.Item(122, 10).Value = 1199.2222
$ws.Cells.Item(122, 12).Value = 10792.9998
$ws.Cells.Item(122, 14).Value = -15692.9998

$ws.Cells.Item(135, 8).Value = 1047.7059
$ws.Cells.Item(135, 9).Value = 482.64285
$ws.Cells.Item(135, 10).Value = 1443.25
$ws.Cells.Item(135, 11).Value = 4343.78565
$ws.Cells.Item(135, 12).Value = 12989.25
$ws.Cells.Item(135, 13).Value = -1808.78565
$ws.Cells.Item(135, 14).Value = -18059.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 0
$ws.Cells.Item(113, 9).Value = 0
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 0
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).ClearContents()
$ws.Cells.Item(113, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 1113951
$ws.Cells.Item(122, 10).Value = 2599.75
$ws.Cells.Item(122, 12).Value = 7799.25
$ws.Cells.Item(122, 14).Value = -12699.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1027
$ws.Cells.Item(16, 9).Value = 1027
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 1027
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -857
$ws.Cells.Item(16, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 3680.0454
$ws.Cells.Item(122, 9).Value = 2327.3333
$ws.Cells.Item(122, 10).Value = 3893.6316
$ws.Cells.Item(122, 11).Value = 6981.999899999999
$ws.Cells.Item(122, 12).Value = 11680.8948
$ws.Cells.Item(122, 13).Value = -4531.999899999999
$ws.Cells.Item(122, 14).Value = -16580.8948

$ws.Cells.Item(132, 8).Value = 4007.3684
$ws.Cells.Item(132, 9).Value = 3541.5908
$ws.Cells.Item(132, 10).Value = 4647.8125
$ws.Cells.Item(132, 11).Value = 10624.7724
$ws.Cells.Item(132, 12).Value = 13943.4375
$ws.Cells.Item(132, 13).Value = -8094.7724
$ws.Cells.Item(132, 14).Value = -19003.4375

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 142678.62
$ws.Cells.Item(46, 10).Value = 142678.62
$ws.Cells.Item(46, 12).Value = 142678.62
$ws.Cells.Item(46, 14).Value = -143140.62

$ws.Cells.Item(132, 8).Value = 13890789
$ws.Cells.Item(132, 9).Value = 20001514
$ws.Cells.Item(132, 10).Value = 2777.9092
$ws.Cells.Item(132, 11).Value = 60004542
$ws.Cells.Item(132, 12).Value = 8333.7276
$ws.Cells.Item(132, 13).Value = -60002012
$ws.Cells.Item(132, 14).Value = -13393.7276

$ws.Cells.Item(134, 8).Value = 142678.62
$ws.Cells.Item(134, 10).Value = 142678.62
$ws.Cells.Item(134, 12).Value = 428035.86
$ws.Cells.Item(134, 14).Value = -433105.86
